$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark from its original location
#    (it gets re-created further below, at the end of the newly typed
#    "." in the approval table, mirroring where Word last left the cursor).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Fill in the "Diagramme de classe" approval row of the tasks table
#    (4th table in the document).
$table = $d.Tables(4)
$row = 6

# Approbateurs column -> "Alban" / "PAPASSIAN" (two separate paragraphs)
$approverCell = $table.Cell($row, 3)
$approverCell.Range.Text = "Alban" + [char]13 + "PAPASSIAN"

# Rendu reel column -> "28/10"
$renduReelCell = $table.Cell($row, 6)
$renduReelCell.Range.Text = "28/10"

# Notes column -> "." followed by a fresh "_GoBack" bookmark right after it.
# A placeholder character is typed after the ".", bookmarked, and then
# removed again - this leaves the bookmark collapsed exactly after the
# "." (matching where Word stashes "_GoBack" after the last edit).
$notesCell = $table.Cell($row, 7)
$notesCell.Range.Text = "." + "X"
$notesRange = $table.Cell($row, 7).Range
$placeholder = $d.Range($notesRange.Start + 1, $notesRange.Start + 2)
$d.Bookmarks.Add("_GoBack", $placeholder)
$toDelete = $d.Range($notesRange.Start + 1, $notesRange.Start + 2)
$toDelete.Text = ""
